$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3500
$ws.Range("I69").Value = 3500
$ws.Range("J69").Value = 3500
$ws.Range("K69").Value = 10500
$ws.Range("L69").Value = 10500
$ws.Range("M69").Value = -9626
$ws.Range("N69").Value = -12248
$ws.Range("H72").Value = 3500
$ws.Range("I72").Value = 3500
$ws.Range("J72").Value = 3500
$ws.Range("K72").Value = 31500
$ws.Range("L72").Value = 31500
$ws.Range("M72").Value = -27132
$ws.Range("N72").Value = -40236
$ws.Range("H106").Value = 266071.84
$ws.Range("I106").Value = 335965.66
$ws.Range("J106").Value = 3970
$ws.Range("K106").Value = 335965.66
$ws.Range("L106").Value = 3970
$ws.Range("M106").Value = -335334.66
$ws.Range("N106").Value = -5232
$ws.Range("H138").Value = 6580829.5
$ws.Range("I138").Value = 1482.125
$ws.Range("K138").Value = 4446.375
$ws.Range("M138").Value = 693.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8270.411
$ws.Range("I32").Value = 9710.880999999999
$ws.Range("J32").Value = 3949
$ws.Range("K32").Value = 9710.880999999999
$ws.Range("L32").Value = 3949
$ws.Range("M32").Value = -9423.880999999999
$ws.Range("N32").Value = -4523
$ws.Range("H61").Value = 17244504
$ws.Range("I61").Value = 17860184
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 17860184
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -17859972
$ws.Range("N61").Value = -5924
$ws.Range("H74").Value = 19234582
$ws.Range("I74").Value = 27779804
$ws.Range("J74").Value = 7830.25
$ws.Range("K74").Value = 27779804
$ws.Range("L74").Value = 7830.25
$ws.Range("M74").Value = -27778930
$ws.Range("N74").Value = -9578.25
$ws.Range("H77").Value = 19234582
$ws.Range("I77").Value = 27779804
$ws.Range("J77").Value = 7830.25
$ws.Range("K77").Value = 138899020
$ws.Range("L77").Value = 39151.25
$ws.Range("M77").Value = -138894652
$ws.Range("N77").Value = -47887.25
$ws.Range("H132").Value = 7145224
$ws.Range("I132").Value = 7814463.5
$ws.Range("J132").Value = 6670.6665
$ws.Range("K132").Value = 23443390.5
$ws.Range("L132").Value = 20011.9995
$ws.Range("M132").Value = -23440860.5
$ws.Range("N132").Value = -25071.9995
$ws.Range("H136").Value = 17244504
$ws.Range("I136").Value = 17860184
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 53580552
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -53578002
$ws.Range("N136").Value = -21600

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 406.62964
$ws.Range("I107").Value = 332.3846
$ws.Range("J107").Value = 475.57144
$ws.Range("K107").Value = 332.3846
$ws.Range("L107").Value = 475.57144
$ws.Range("M107").Value = 1587.6154
$ws.Range("N107").Value = -4315.57144
$ws.Range("H134").Value = 1136036.4
$ws.Range("I134").Value = 2095.6428
$ws.Range("J134").Value = 3403917.8
$ws.Range("K134").Value = 6286.928400000001
$ws.Range("L134").Value = 10211753.4
$ws.Range("M134").Value = -3751.928400000001
$ws.Range("N134").Value = -10216823.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1576
$ws.Range("I114").Value = 350
$ws.Range("J114").Value = 1882.5
$ws.Range("K114").Value = 1050
$ws.Range("L114").Value = 5647.5
$ws.Range("M114").Value = 2204
$ws.Range("N114").Value = -12155.5
$ws.Range("H131").Value = 2398887.8
$ws.Range("I131").Value = 9592759
$ws.Range("J131").Value = 930.61536
$ws.Range("K131").Value = 28778277
$ws.Range("L131").Value = 2791.84608
$ws.Range("M131").Value = -28773237
$ws.Range("N131").Value = -12871.84608

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19426190
$ws.Range("I80").Value = 23812268
$ws.Range("J80").Value = 7145177
$ws.Range("K80").Value = 23812268
$ws.Range("L80").Value = 7145177
$ws.Range("M80").Value = -23811270
$ws.Range("N80").Value = -7147173
$ws.Range("H83").Value = 19426190
$ws.Range("I83").Value = 23812268
$ws.Range("J83").Value = 7145177
$ws.Range("K83").Value = 119061340
$ws.Range("L83").Value = 35725885
$ws.Range("M83").Value = -119056348
$ws.Range("N83").Value = -35735869
$ws.Range("H113").Value = 167585.5
$ws.Range("I113").Value = 251015.25
$ws.Range("J113").Value = 726
$ws.Range("K113").Value = 251015.25
$ws.Range("L113").Value = 726
$ws.Range("M113").Value = -248845.25
$ws.Range("N113").Value = -5066
$ws.Range("H126").Value = 3236.0571
$ws.Range("I126").Value = 1877.0526
$ws.Range("J126").Value = 4849.875
$ws.Range("K126").Value = 5631.1578
$ws.Range("L126").Value = 14549.625
$ws.Range("M126").Value = -3161.1578
$ws.Range("N126").Value = -19489.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 946.8
$ws.Range("I46").Value = 955.5
$ws.Range("J46").Value = 945.46155
$ws.Range("K46").Value = 955.5
$ws.Range("L46").Value = 945.46155
$ws.Range("M46").Value = -767.5
$ws.Range("N46").Value = -1321.46155
$ws.Range("H61").Value = 1618
$ws.Range("I61").Value = 1476.1538
$ws.Range("J61").Value = 1822.8889
$ws.Range("K61").Value = 1476.1538
$ws.Range("L61").Value = 1822.8889
$ws.Range("M61").Value = -1274.1538
$ws.Range("N61").Value = -2226.8889
$ws.Range("H68").Value = 1891.7368
$ws.Range("J68").Value = 2185.3635
$ws.Range("L68").Value = 2185.3635
$ws.Range("N68").Value = -3683.3635
$ws.Range("H71").Value = 1891.7368
$ws.Range("J71").Value = 2185.3635
$ws.Range("L71").Value = 10926.8175
$ws.Range("N71").Value = -18414.8175
$ws.Range("H100").Value = 1907.0625
$ws.Range("I100").Value = 1756.5555
$ws.Range("J100").Value = 2100.5715
$ws.Range("K100").Value = 1756.5555
$ws.Range("L100").Value = 2100.5715
$ws.Range("M100").Value = -1215.5555
$ws.Range("N100").Value = -3182.5715
$ws.Range("H113").Value = 1618
$ws.Range("I113").Value = 1476.1538
$ws.Range("J113").Value = 1822.8889
$ws.Range("K113").Value = 1476.1538
$ws.Range("L113").Value = 1822.8889
$ws.Range("M113").Value = 693.8462
$ws.Range("N113").Value = -6162.8889

# ---- Sheet: WVR (special rows with added/removed cells) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1000
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1284

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
